# Updates "想去人数" (want-to-go count) figures across all four sheets of the
# 广州-漫展信息 workbook, matching the refreshed scrape output at commit 456a3b4.
#
# Sheet order (tab order == workbook.xml <sheets> order):
#   1 = 展览     (Exhibitions)
#   2 = 演出     (Performances)
#   3 = 本地生活 (Local life)
#   4 = 全部类型 (All types - date-sorted union of the first three sheets)

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 (Exhibitions) ----
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 8290
$ws.Range("F3").Value = 130
$ws.Range("F4").Value = 102
$ws.Range("F5").Value = 35624
$ws.Range("F7").Value = 611
$ws.Range("F8").Value = 725
$ws.Range("F10").Value = 150
$ws.Range("F12").Value = 815
$ws.Range("F13").Value = 68
$ws.Range("F14").Value = 637
$ws.Range("F15").Value = 447
$ws.Range("F17").Value = 584
$ws.Range("F18").Value = 161
$ws.Range("F19").Value = 432
$ws.Range("F20").Value = 427
$ws.Range("F21").Value = 1124
$ws.Range("F23").Value = 743
$ws.Range("F24").Value = 2398
$ws.Range("F25").Value = 880
$ws.Range("F26").Value = 508
$ws.Range("F27").Value = 78
$ws.Range("F28").Value = 1103
$ws.Range("F30").Value = 667
$ws.Range("F31").Value = 667
$ws.Range("F33").Value = 1102

# ---- Sheet 2: 演出 (Performances) ----
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 300

# ---- Sheet 3: 本地生活 (Local life) ----
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 557

# ---- Sheet 4: 全部类型 (All types, union of sheets 1-3, date sorted) ----
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 557
$ws.Range("F3").Value = 8290
$ws.Range("F4").Value = 130
$ws.Range("F5").Value = 102
$ws.Range("F6").Value = 300
$ws.Range("F7").Value = 35624
# COMICUP 2024SP (row 7) is no longer sold out - now has a real min price.
$ws.Range("G7").Value = 68
$ws.Range("F9").Value = 611
$ws.Range("F10").Value = 725
$ws.Range("F13").Value = 150
$ws.Range("F18").Value = 815
$ws.Range("F19").Value = 68
$ws.Range("F20").Value = 637
$ws.Range("F21").Value = 447
$ws.Range("F28").Value = 584
$ws.Range("F29").Value = 161
$ws.Range("F30").Value = 432
$ws.Range("F31").Value = 427
$ws.Range("F32").Value = 1124
$ws.Range("F34").Value = 743
$ws.Range("F35").Value = 2398
$ws.Range("F36").Value = 880
$ws.Range("F37").Value = 508
$ws.Range("F38").Value = 78
$ws.Range("F39").Value = 1103
$ws.Range("F42").Value = 667
$ws.Range("F43").Value = 667
$ws.Range("F45").Value = 1102
